$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 20:20"

# Row 4
$ws.Range("B4").Value = 155969
$ws.Range("C4").Value = 12478
$ws.Range("E4").Value = 147904
$ws.Range("G4").Value = 271
$ws.Range("H4").Value = 2854

# Row 9
$ws.Range("A9").Value = "Francia"
$ws.Range("B9").Value = 44550
$ws.Range("C9").Value = 4376
$ws.Range("D9").Value = 7927
$ws.Range("E9").Value = 33599
$ws.Range("F9").Value = 5056
$ws.Range("G9").Value = 418
$ws.Range("H9").Value = 3024

# Row 10
$ws.Range("A10").Value = "Iran"
$ws.Range("B10").Value = 41495
$ws.Range("C10").Value = 3186
$ws.Range("D10").Value = 13911
$ws.Range("E10").Value = 24827
$ws.Range("F10").Value = 3511
$ws.Range("G10").Value = 117
$ws.Range("H10").Value = 2757

# Row 30
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 2055
$ws.Range("C30").Value = 193
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = 2017
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = 31

# Row 31
$ws.Range("A31").Value = "Luxemburgo"
$ws.Range("B31").Value = 1988
$ws.Range("C31").Value = 38
$ws.Range("D31").Value = 40
$ws.Range("E31").Value = 1926
$ws.Range("F31").Value = 25
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 22

# Row 32
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 1962
$ws.Range("C32").Value = 38
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 1899
$ws.Range("F32").Value = 58
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 60

# Row 33
$ws.Range("A33").Value = "Rumania"
$ws.Range("B33").Value = 1952
$ws.Range("C33").Value = 137
$ws.Range("D33").Value = 209
$ws.Range("E33").Value = 1691
$ws.Range("F33").Value = 33
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 52

# Row 48
$ws.Range("A48").Value = "Peru"
$ws.Range("B48").Value = 950
$ws.Range("C48").Value = 98
$ws.Range("D48").Value = 16
$ws.Range("E48").Value = 916
$ws.Range("F48").Value = 40
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 18

# Row 49
$ws.Range("A49").Value = "Republica Dominicana"
$ws.Range("B49").Value = 901
$ws.Range("C49").Value = 42
$ws.Range("D49").Value = 4
$ws.Range("E49").Value = 855
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 42

# Row 50
$ws.Range("A50").Value = "Singapur"
$ws.Range("B50").Value = 879
$ws.Range("C50").Value = 35
$ws.Range("D50").Value = 228
$ws.Range("E50").Value = 648
$ws.Range("F50").Value = 19
$ws.Range("H50").Value = 3

# Row 51
$ws.Range("E51").Value = 569
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 23

# Row 59
$ws.Range("A59").Value = "Egipto"
$ws.Range("B59").Value = 656
$ws.Range("C59").Value = 47
$ws.Range("D59").Value = 150
$ws.Range("E59").Value = 465
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 41

# Row 60
$ws.Range("A60").Value = "Hong Kong"
$ws.Range("B60").Value = 642
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 118
$ws.Range("E60").Value = 520
$ws.Range("F60").Value = 5
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 4

# Row 61
$ws.Range("A61").Value = "Irak"
$ws.Range("B61").Value = 630
$ws.Range("C61").Value = 83
$ws.Range("D61").Value = 152
$ws.Range("E61").Value = 432
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 46

# Row 62
$ws.Range("A62").Value = "Emiratos Arabes Unidos"
$ws.Range("B62").Value = 611
$ws.Range("C62").Value = 41
$ws.Range("D62").Value = 61
$ws.Range("E62").Value = 545
$ws.Range("F62").Value = 2
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 5

# Row 67
$ws.Range("D67").Value = 7
$ws.Range("E67").Value = 477

# Row 88
$ws.Range("D88").Value = 22
$ws.Range("E88").Value = 201

# Row 113
$ws.Range("B113").Value = 103
$ws.Range("C113").Value = 12
$ws.Range("D113").Value = 20
$ws.Range("E113").Value = 83

# Row 130
$ws.Range("A130").Value = "Madagascar"
$ws.Range("C130").Value = 4
$ws.Range("D130").Value = 0
$ws.Range("E130").Value = 43

# Row 131
$ws.Range("A131").Value = "Guayana Francesa"
$ws.Range("B131").Value = 43
$ws.Range("D131").Value = 6
$ws.Range("E131").Value = 37

# Row 136
$ws.Range("A136").Value = "Zambia"
$ws.Range("C136").Value = 6
$ws.Range("F136").Value = 0

# Row 137
$ws.Range("A137").Value = "Polinesia Francesa"
$ws.Range("C137").Value = 5
$ws.Range("F137").Value = 2

# Row 138
$ws.Range("A138").Value = "Barbados"

# Row 139
$ws.Range("A139").Value = "Uganda"

# Row 168
$ws.Range("A168").Value = "Siria"
$ws.Range("D168").Value = 0
$ws.Range("H168").Value = 1

# Row 169
$ws.Range("A169").Value = "Santa Lucia"
$ws.Range("D169").Value = 1
$ws.Range("H169").Value = 0

# Row 170
$ws.Range("A170").Value = "Libia"

# Row 171
$ws.Range("A171").Value = "Mozambique"

# Row 172
$ws.Range("A172").Value = "Seychelles"

# Row 173
$ws.Range("A173").Value = "Guinea-Bisau"
$ws.Range("C173").Value = 6

# Row 174
$ws.Range("A174").Value = "Laos"

# Row 175
$ws.Range("A175").Value = "Surinam"
$ws.Range("C175").Value = 0

# Row 182
$ws.Range("A182").Value = "Santa Sede"

# Row 183
$ws.Range("A183").Value = "San Martin (Parte Holandesa)"

# Row 184
$ws.Range("A184").Value = "Benin"

# Row 185
$ws.Range("A185").Value = "Cabo Verde"
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("H185").Value = 1

# Row 186
$ws.Range("A186").Value = "San Bartolome"
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 1
$ws.Range("H186").Value = 0

# Row 188
$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("C188").Value = 1

# Row 190
$ws.Range("A190").Value = "Montserrat"
$ws.Range("C190").Value = 0

# Row 191
$ws.Range("A191").Value = "Republica del Chad"
$ws.Range("C191").Value = 2

# Row 197
$ws.Range("A197").Value = "Botsuana"
$ws.Range("C197").Value = 3

# Row 198
$ws.Range("A198").Value = "Belice"
$ws.Range("C198").Value = 1

# Row 201
$ws.Range("A201").Value = "Somalia"
$ws.Range("B201").Value = 3
$ws.Range("E201").Value = 3

# Row 203
$ws.Range("A203").Value = "Anguila"
$ws.Range("B203").Value = 2
$ws.Range("E203").Value = 2

# Row 205
$ws.Range("A205").Value = "Papua Nueva Guinea"
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 1

# Row 206
$ws.Range("A206").Value = "San Vicente y las Granadinas"
$ws.Range("B206").Value = 1
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 1
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
